# ---------------------------------------------------------------------------
# "updated figure and data"
#
# The old "ReFuelEU" sheet (outdated SAF targets) is kept but renamed to
# "ReFuelEU (outdated)". A brand-new "ReFuelEU" sheet with refreshed figures
# (and fresh EUR-Lex source links) is inserted in front of it. "EU Production"
# stays last and loses tab-selection in favour of the new first sheet.
#
# NOTE: worksheet variables in this host resolve by tab *position*, not by a
# stable object handle -- once Worksheets.Add()/Move() shifts tab order, an
# old variable silently starts pointing at whatever sheet now sits at that
# position. To stay correct we always re-fetch sheets by name right before
# using them instead of reusing a variable captured earlier.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# 1) rename the existing (now superseded) sheet
$wb.Worksheets.Item("ReFuelEU").Name = "ReFuelEU (outdated)"

# 2) insert the new sheet in front of it -> becomes the first tab
$wsNew = $wb.Worksheets.Add($wb.Worksheets.Item("ReFuelEU (outdated)"))
$wsNew.Name = "ReFuelEU"

# --- populate the new "ReFuelEU" sheet -------------------------------------

$wb.Worksheets.Item("ReFuelEU").Range("A1").Value = "year"
$wb.Worksheets.Item("ReFuelEU").Range("B1").Value = "SAF share [%]"
$wb.Worksheets.Item("ReFuelEU").Range("C1").Value = "of which synth-fuel share [%]"
$wb.Worksheets.Item("ReFuelEU").Range("D1").Value = "source"

$wb.Worksheets.Item("ReFuelEU").Range("A2").Value = 2025
$wb.Worksheets.Item("ReFuelEU").Range("B2").Value = 2
$wb.Worksheets.Item("ReFuelEU").Range("C2").Value = 0
$wb.Worksheets.Item("ReFuelEU").Range("D2").Value = "https://eur-lex.europa.eu/legal-content/EN/TXT/HTML/?uri=OJ:L_202302405"

$wb.Worksheets.Item("ReFuelEU").Range("A3").Value = 2030
$wb.Worksheets.Item("ReFuelEU").Range("B3").Value = 6
$wb.Worksheets.Item("ReFuelEU").Range("C3").Value = 0.7
$wb.Worksheets.Item("ReFuelEU").Range("D3").Value = "https://eur-lex.europa.eu/legal-content/EN/TXT/HTML/?uri=OJ:L_202302406"

$wb.Worksheets.Item("ReFuelEU").Range("A4").Value = 2035
$wb.Worksheets.Item("ReFuelEU").Range("B4").Value = 20
$wb.Worksheets.Item("ReFuelEU").Range("C4").Value = 5
$wb.Worksheets.Item("ReFuelEU").Range("D4").Value = "https://eur-lex.europa.eu/legal-content/EN/TXT/HTML/?uri=OJ:L_202302407"

$wb.Worksheets.Item("ReFuelEU").Range("A5").Value = 2040
$wb.Worksheets.Item("ReFuelEU").Range("B5").Value = 34
$wb.Worksheets.Item("ReFuelEU").Range("C5").Value = 10
$wb.Worksheets.Item("ReFuelEU").Range("D5").Value = "https://eur-lex.europa.eu/legal-content/EN/TXT/HTML/?uri=OJ:L_202302408"

$wb.Worksheets.Item("ReFuelEU").Range("A6").Value = 2045
$wb.Worksheets.Item("ReFuelEU").Range("B6").Value = 42
$wb.Worksheets.Item("ReFuelEU").Range("C6").Value = 15
$wb.Worksheets.Item("ReFuelEU").Range("D6").Value = "https://eur-lex.europa.eu/legal-content/EN/TXT/HTML/?uri=OJ:L_202302409"

$wb.Worksheets.Item("ReFuelEU").Range("A7").Value = 2050
$wb.Worksheets.Item("ReFuelEU").Range("B7").Value = 70
$wb.Worksheets.Item("ReFuelEU").Range("C7").Value = 35
$wb.Worksheets.Item("ReFuelEU").Range("D7").Value = "https://eur-lex.europa.eu/legal-content/EN/TXT/HTML/?uri=OJ:L_202302410"

# style the source column like the old sheet's hyperlink-look cells
$wb.Worksheets.Item("ReFuelEU").Range("D2:D7").Style = "Hyperlink"

# column widths matching the "outdated" sheet's layout
$wb.Worksheets.Item("ReFuelEU").Columns.Item(2).ColumnWidth = 11.830729166666666
$wb.Worksheets.Item("ReFuelEU").Columns.Item(3).ColumnWidth = 24.666666666666668

# --- tidy up the other two sheets' view state -------------------------------

# "ReFuelEU (outdated)" keeps its data but gains the same column widths and a
# refreshed selection; it is no longer the active tab.
$wb.Worksheets.Item("ReFuelEU (outdated)").Columns.Item(2).ColumnWidth = 11.830729166666666
$wb.Worksheets.Item("ReFuelEU (outdated)").Columns.Item(3).ColumnWidth = 24.666666666666668
$wb.Worksheets.Item("ReFuelEU (outdated)").Range("D13").Select()

# "EU Production" no longer the active tab either; selection unchanged.
$wb.Worksheets.Item("EU Production").Range("E42").Select()

# view state: selection + zoom, then make this the tab shown on open -- do
# this LAST since selecting a range on a sheet also activates that sheet.
$wb.Worksheets.Item("ReFuelEU").Range("C10").Select()
$wb.Worksheets.Item("ReFuelEU").Activate()
$excel.ActiveWindow.Zoom = 135
